$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: add the new commit entry "equip instruction & functionality added"
# with 2 hours logged (this becomes shared-string index 18).
$ws.Range("C18").Value = "equip instruction & functionality added"
$ws.Range("G18").Value = 2

# Match the commit-row look used by the other entries (C5:C17), i.e. the
# "20% - Accent5" cell style already applied throughout that column.
$ws.Range("C18").Style = $ws.Range("C17").Style

# Move the active selection to H18, as recorded in the saved view state.
$ws.Range("H18").Select()
